$wb = $excel.ActiveWorkbook

# --- Sheet "Recommandations": update rows 2-31, then delete old rows 32-42 ---
$ws1 = $wb.Worksheets.Item("Recommandations")

# Row 2: BRVM - CONSOMMATION DISCRETIONNAIRE
$ws1.Cells.Item(2,1).Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws1.Cells.Item(2,2).Value = 0
$ws1.Cells.Item(2,3).Value = 3
$ws1.Cells.Item(2,4).Value = 492.58
$ws1.Cells.Item(2,5).Value = 158.97
$ws1.Cells.Item(2,6).Value = "🟡 Observer"
$ws1.Cells.Item(2,7).Value = "➖ Neutre"

# Row 3: BRVM - SERVICES FINANCIERS
$ws1.Cells.Item(3,1).Value = "BRVM - SERVICES FINANCIERS"
$ws1.Cells.Item(3,2).Value = 0
$ws1.Cells.Item(3,3).Value = 3
$ws1.Cells.Item(3,4).Value = 436.3
$ws1.Cells.Item(3,5).Value = 145.69
$ws1.Cells.Item(3,6).Value = "🟡 Observer"
$ws1.Cells.Item(3,7).Value = "➖ Neutre"

# Row 4: BRVM - CONSOMMATION DE BASE    (**)
$ws1.Cells.Item(4,1).Value = "BRVM - CONSOMMATION DE BASE    (**)"
$ws1.Cells.Item(4,2).Value = 0
$ws1.Cells.Item(4,3).Value = 2
$ws1.Cells.Item(4,4).Value = 431.6
$ws1.Cells.Item(4,5).Value = 215.68
$ws1.Cells.Item(4,6).Value = "🟡 Observer"
$ws1.Cells.Item(4,7).Value = "➖ Neutre"

# Row 5: BRVM-PRESTIGE
$ws1.Cells.Item(5,1).Value = "BRVM-PRESTIGE"
$ws1.Cells.Item(5,2).Value = 0
$ws1.Cells.Item(5,3).Value = 3
$ws1.Cells.Item(5,4).Value = 428.05
$ws1.Cells.Item(5,5).Value = 143.16
$ws1.Cells.Item(5,6).Value = "🟡 Observer"
$ws1.Cells.Item(5,7).Value = "➖ Neutre"

# Row 6: BRVM - INDUSTRIELS
$ws1.Cells.Item(6,1).Value = "BRVM - INDUSTRIELS"
$ws1.Cells.Item(6,2).Value = 0
$ws1.Cells.Item(6,3).Value = 3
$ws1.Cells.Item(6,4).Value = 393.69
$ws1.Cells.Item(6,5).Value = 130.52
$ws1.Cells.Item(6,6).Value = "🟡 Observer"
$ws1.Cells.Item(6,7).Value = "➖ Neutre"

# Row 7: BRVM - ENERGIE
$ws1.Cells.Item(7,1).Value = "BRVM - ENERGIE"
$ws1.Cells.Item(7,2).Value = 0
$ws1.Cells.Item(7,3).Value = 3
$ws1.Cells.Item(7,4).Value = 331.04
$ws1.Cells.Item(7,5).Value = 111.1
$ws1.Cells.Item(7,6).Value = "🟡 Observer"
$ws1.Cells.Item(7,7).Value = "➖ Neutre"

# Row 8: BRVM - SERVICES PUBLICS
$ws1.Cells.Item(8,1).Value = "BRVM - SERVICES PUBLICS"
$ws1.Cells.Item(8,2).Value = 0
$ws1.Cells.Item(8,3).Value = 3
$ws1.Cells.Item(8,4).Value = 319.72
$ws1.Cells.Item(8,5).Value = 106.38
$ws1.Cells.Item(8,6).Value = "🟡 Observer"
$ws1.Cells.Item(8,7).Value = "➖ Neutre"

# Row 9: BRVM - TELECOMMUNICATIONS
$ws1.Cells.Item(9,1).Value = "BRVM - TELECOMMUNICATIONS"
$ws1.Cells.Item(9,2).Value = 0
$ws1.Cells.Item(9,3).Value = 3
$ws1.Cells.Item(9,4).Value = 279.61
$ws1.Cells.Item(9,5).Value = 93.5
$ws1.Cells.Item(9,6).Value = "🟡 Observer"
$ws1.Cells.Item(9,7).Value = "➖ Neutre"

# Row 10: BRVM – COMPOSITE TOTAL RETURN    (**)
$ws1.Cells.Item(10,1).Value = "BRVM – COMPOSITE TOTAL RETURN    (**)"
$ws1.Cells.Item(10,2).Value = 0
$ws1.Cells.Item(10,3).Value = 2
$ws1.Cells.Item(10,4).Value = 263.52
$ws1.Cells.Item(10,5).Value = 130.87
$ws1.Cells.Item(10,6).Value = "🟡 Observer"
$ws1.Cells.Item(10,7).Value = "➖ Neutre"

# Row 11: BRVM-PRINCIPAL    (**)
$ws1.Cells.Item(11,1).Value = "BRVM-PRINCIPAL    (**)"
$ws1.Cells.Item(11,2).Value = 0
$ws1.Cells.Item(11,3).Value = 1
$ws1.Cells.Item(11,4).Value = 217.75
$ws1.Cells.Item(11,5).Value = 217.75
$ws1.Cells.Item(11,6).Value = "🟡 Observer"
$ws1.Cells.Item(11,7).Value = "➖ Neutre"

# Row 12: BRVM-PRINCIPAL     (**)
$ws1.Cells.Item(12,1).Value = "BRVM-PRINCIPAL     (**)"
$ws1.Cells.Item(12,2).Value = 0
$ws1.Cells.Item(12,3).Value = 1
$ws1.Cells.Item(12,4).Value = 216.13
$ws1.Cells.Item(12,5).Value = 216.13
$ws1.Cells.Item(12,6).Value = "🟡 Observer"
$ws1.Cells.Item(12,7).Value = "➖ Neutre"

# Row 13: EVIOSYS PACKAGING SIEM CI (SEMC)
$ws1.Cells.Item(13,1).Value = "EVIOSYS PACKAGING SIEM CI (SEMC)"
$ws1.Cells.Item(13,2).Value = 3
$ws1.Cells.Item(13,3).Value = 0
$ws1.Cells.Item(13,4).Value = 21.92
$ws1.Cells.Item(13,5).Value = 7.45
$ws1.Cells.Item(13,6).Value = "🟢 Achat"
$ws1.Cells.Item(13,7).Value = "✅ Renforcer"

# Row 14: ORAGROUP TOGO (ORGT)
$ws1.Cells.Item(14,1).Value = "ORAGROUP TOGO (ORGT)"
$ws1.Cells.Item(14,2).Value = 2
$ws1.Cells.Item(14,3).Value = 0
$ws1.Cells.Item(14,4).Value = 9.58
$ws1.Cells.Item(14,5).Value = 6.04
$ws1.Cells.Item(14,6).Value = "🟡 Observer"
$ws1.Cells.Item(14,7).Value = "➖ Neutre"

# Row 15: ERIUM CI (Ex AIR LIQUIDE CI) (SIVC)
$ws1.Cells.Item(15,1).Value = "ERIUM CI (Ex AIR LIQUIDE CI) (SIVC)"
$ws1.Cells.Item(15,2).Value = 2
$ws1.Cells.Item(15,3).Value = 0
$ws1.Cells.Item(15,4).Value = 8.41
$ws1.Cells.Item(15,5).Value = 2.08
$ws1.Cells.Item(15,6).Value = "🟡 Observer"
$ws1.Cells.Item(15,7).Value = "➖ Neutre"

# Row 16: SICOR CI (SICC)
$ws1.Cells.Item(16,1).Value = "SICOR CI (SICC)"
$ws1.Cells.Item(16,2).Value = 1
$ws1.Cells.Item(16,3).Value = 0
$ws1.Cells.Item(16,4).Value = 7.4
$ws1.Cells.Item(16,5).Value = 7.4
$ws1.Cells.Item(16,6).Value = "🟡 Observer"
$ws1.Cells.Item(16,7).Value = "➖ Neutre"

# Row 17: BANK OF AFRICA BF (BOABF)
$ws1.Cells.Item(17,1).Value = "BANK OF AFRICA BF (BOABF)"
$ws1.Cells.Item(17,2).Value = 1
$ws1.Cells.Item(17,3).Value = 0
$ws1.Cells.Item(17,4).Value = 5.07
$ws1.Cells.Item(17,5).Value = 5.07
$ws1.Cells.Item(17,6).Value = "🟡 Observer"
$ws1.Cells.Item(17,7).Value = "➖ Neutre"

# Row 18: VIVO ENERGY CI (SHEC)
$ws1.Cells.Item(18,1).Value = "VIVO ENERGY CI (SHEC)"
$ws1.Cells.Item(18,2).Value = 1
$ws1.Cells.Item(18,3).Value = 0
$ws1.Cells.Item(18,4).Value = 4.96
$ws1.Cells.Item(18,5).Value = 4.96
$ws1.Cells.Item(18,6).Value = "🟡 Observer"
$ws1.Cells.Item(18,7).Value = "➖ Neutre"

# Row 19: SOGB CI (SOGC)
$ws1.Cells.Item(19,1).Value = "SOGB CI (SOGC)"
$ws1.Cells.Item(19,2).Value = 1
$ws1.Cells.Item(19,3).Value = 0
$ws1.Cells.Item(19,4).Value = 3.24
$ws1.Cells.Item(19,5).Value = 3.24
$ws1.Cells.Item(19,6).Value = "🟡 Observer"
$ws1.Cells.Item(19,7).Value = "➖ Neutre"

# Row 20: NESTLE CI (NTLC)
$ws1.Cells.Item(20,1).Value = "NESTLE CI (NTLC)"
$ws1.Cells.Item(20,2).Value = 1
$ws1.Cells.Item(20,3).Value = 0
$ws1.Cells.Item(20,4).Value = 3.24
$ws1.Cells.Item(20,5).Value = 3.24
$ws1.Cells.Item(20,6).Value = "🟡 Observer"
$ws1.Cells.Item(20,7).Value = "➖ Neutre"

# Row 21: ONATEL BF (ONTBF)
$ws1.Cells.Item(21,1).Value = "ONATEL BF (ONTBF)"
$ws1.Cells.Item(21,2).Value = 1
$ws1.Cells.Item(21,3).Value = 0
$ws1.Cells.Item(21,4).Value = 1.22
$ws1.Cells.Item(21,5).Value = 1.22
$ws1.Cells.Item(21,6).Value = "🟡 Observer"
$ws1.Cells.Item(21,7).Value = "➖ Neutre"

# Row 22: UNIWAX CI (UNXC)
$ws1.Cells.Item(22,1).Value = "UNIWAX CI (UNXC)"
$ws1.Cells.Item(22,2).Value = 0
$ws1.Cells.Item(22,3).Value = 1
$ws1.Cells.Item(22,4).Value = -2.13
$ws1.Cells.Item(22,5).Value = -2.13
$ws1.Cells.Item(22,6).Value = "🟡 Observer"
$ws1.Cells.Item(22,7).Value = "➖ Neutre"

# Row 23: BERNABE CI (BNBC)
$ws1.Cells.Item(23,1).Value = "BERNABE CI (BNBC)"
$ws1.Cells.Item(23,2).Value = 0
$ws1.Cells.Item(23,3).Value = 1
$ws1.Cells.Item(23,4).Value = -3.19
$ws1.Cells.Item(23,5).Value = -3.19
$ws1.Cells.Item(23,6).Value = "🟡 Observer"
$ws1.Cells.Item(23,7).Value = "➖ Neutre"

# Row 24: SETAO CI (STAC)
$ws1.Cells.Item(24,1).Value = "SETAO CI (STAC)"
$ws1.Cells.Item(24,2).Value = 0
$ws1.Cells.Item(24,3).Value = 1
$ws1.Cells.Item(24,4).Value = -3.4
$ws1.Cells.Item(24,5).Value = -3.4
$ws1.Cells.Item(24,6).Value = "🟡 Observer"
$ws1.Cells.Item(24,7).Value = "➖ Neutre"

# Row 25: SOLIBRA CI (SLBC)
$ws1.Cells.Item(25,1).Value = "SOLIBRA CI (SLBC)"
$ws1.Cells.Item(25,2).Value = 1
$ws1.Cells.Item(25,3).Value = 1
$ws1.Cells.Item(25,4).Value = -3.45
$ws1.Cells.Item(25,5).Value = -4.5
$ws1.Cells.Item(25,6).Value = "🟡 Observer"
$ws1.Cells.Item(25,7).Value = "👀 À surveiller"

# Row 26: TOTALENERGIES MARKETING SN (TTLS)
$ws1.Cells.Item(26,1).Value = "TOTALENERGIES MARKETING SN (TTLS)"
$ws1.Cells.Item(26,2).Value = 1
$ws1.Cells.Item(26,3).Value = 1
$ws1.Cells.Item(26,4).Value = -3.92
$ws1.Cells.Item(26,5).Value = 3.19
$ws1.Cells.Item(26,6).Value = "🟡 Observer"
$ws1.Cells.Item(26,7).Value = "👀 À surveiller"

# Row 27: SAFCA CI (SAFC)
$ws1.Cells.Item(27,1).Value = "SAFCA CI (SAFC)"
$ws1.Cells.Item(27,2).Value = 0
$ws1.Cells.Item(27,3).Value = 1
$ws1.Cells.Item(27,4).Value = -3.93
$ws1.Cells.Item(27,5).Value = -3.93
$ws1.Cells.Item(27,6).Value = "🟡 Observer"
$ws1.Cells.Item(27,7).Value = "➖ Neutre"

# Row 28: ECOBANK TRANS. INCORP. TG (ETIT)
$ws1.Cells.Item(28,1).Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$ws1.Cells.Item(28,2).Value = 0
$ws1.Cells.Item(28,3).Value = 1
$ws1.Cells.Item(28,4).Value = -4.35
$ws1.Cells.Item(28,5).Value = -4.35
$ws1.Cells.Item(28,6).Value = "🟡 Observer"
$ws1.Cells.Item(28,7).Value = "➖ Neutre"

# Row 29: FILTISAC CI (FTSC)
$ws1.Cells.Item(29,1).Value = "FILTISAC CI (FTSC)"
$ws1.Cells.Item(29,2).Value = 0
$ws1.Cells.Item(29,3).Value = 2
$ws1.Cells.Item(29,4).Value = -12.69
$ws1.Cells.Item(29,5).Value = -7.45
$ws1.Cells.Item(29,6).Value = "🟡 Observer"
$ws1.Cells.Item(29,7).Value = "➖ Neutre"

# Row 30: CFAO MOTORS CI (CFAC)
$ws1.Cells.Item(30,1).Value = "CFAO MOTORS CI (CFAC)"
$ws1.Cells.Item(30,2).Value = 0
$ws1.Cells.Item(30,3).Value = 3
$ws1.Cells.Item(30,4).Value = -13.11
$ws1.Cells.Item(30,5).Value = -4.21
$ws1.Cells.Item(30,6).Value = "🔴 Vente"
$ws1.Cells.Item(30,7).Value = "⚠️ Risque de décrochage"

# Row 31: NEI-CEDA CI (NEIC)
$ws1.Cells.Item(31,1).Value = "NEI-CEDA CI (NEIC)"
$ws1.Cells.Item(31,2).Value = 0
$ws1.Cells.Item(31,3).Value = 3
$ws1.Cells.Item(31,4).Value = -17.62
$ws1.Cells.Item(31,5).Value = -2.91
$ws1.Cells.Item(31,6).Value = "🔴 Vente"
$ws1.Cells.Item(31,7).Value = "⚠️ Risque de décrochage"

# Delete now-obsolete trailing rows (old rows 32-42), bottom-up so indices stay valid
$ws1.Rows.Item(42).Delete()
$ws1.Rows.Item(41).Delete()
$ws1.Rows.Item(40).Delete()
$ws1.Rows.Item(39).Delete()
$ws1.Rows.Item(38).Delete()
$ws1.Rows.Item(37).Delete()
$ws1.Rows.Item(36).Delete()
$ws1.Rows.Item(35).Delete()
$ws1.Rows.Item(34).Delete()
$ws1.Rows.Item(33).Delete()
$ws1.Rows.Item(32).Delete()

# --- Sheet "Top_YTD": update rows 2-11 (row count unchanged) ---
$ws2 = $wb.Worksheets.Item("Top_YTD")

# Row 2: BRVM - CONSOMMATION DISCRETIONNAIRE
$ws2.Cells.Item(2,1).Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws2.Cells.Item(2,2).Value = 1743.16

# Row 3: BRVM - SERVICES FINANCIERS
$ws2.Cells.Item(3,1).Value = "BRVM - SERVICES FINANCIERS"
$ws2.Cells.Item(3,2).Value = 1378.42

# Row 4: BRVM-PRESTIGE
$ws2.Cells.Item(4,1).Value = "BRVM-PRESTIGE"
$ws2.Cells.Item(4,2).Value = 1329.26

# Row 5: BRVM - INDUSTRIELS
$ws2.Cells.Item(5,1).Value = "BRVM - INDUSTRIELS"
$ws2.Cells.Item(5,2).Value = 1136.3

# Row 6: BRVM - CONSOMMATION DE BASE    (**)
$ws2.Cells.Item(6,1).Value = "BRVM - CONSOMMATION DE BASE    (**)"
$ws2.Cells.Item(6,2).Value = 897.3

# Row 7: BRVM - ENERGIE
$ws2.Cells.Item(7,1).Value = "BRVM - ENERGIE"
$ws2.Cells.Item(7,2).Value = 830.66

# Row 8: BRVM - SERVICES PUBLICS
$ws2.Cells.Item(8,1).Value = "BRVM - SERVICES PUBLICS"
$ws2.Cells.Item(8,2).Value = 781.49

# Row 9: BRVM - TELECOMMUNICATIONS
$ws2.Cells.Item(9,1).Value = "BRVM - TELECOMMUNICATIONS"
$ws2.Cells.Item(9,2).Value = 621.16

# Row 10: BRVM – COMPOSITE TOTAL RETURN    (**)
$ws2.Cells.Item(10,1).Value = "BRVM – COMPOSITE TOTAL RETURN    (**)"
$ws2.Cells.Item(10,2).Value = 437.12

# Row 11: BRVM-PRINCIPAL    (**)
$ws2.Cells.Item(11,1).Value = "BRVM-PRINCIPAL    (**)"
$ws2.Cells.Item(11,2).Value = 217.75

